$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet used a 2-row header: row 1 held a few (clipped) labels
# and row 2 held unit labels ((m3/s), (MW), (GWh)) spread across F:K.
# The new layout collapses this into a single header row with 11 columns
# (idx, idx2, Name, Date Start, Date End, then the unit/value headers) and
# the data rows shift up by one.

# Remove the old "units" row - everything below shifts up one row and the
# data rows (old 3..14) become the new rows 2..13.
$ws.Rows.Item(2).Delete()

# Start the new header row from a clean slate (the old row 1 left some
# cells with inherited formatting).
$ws.Range("A1:K1").ClearFormats()

# New consolidated header row.
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# F1:K1 use the sheet's secondary "9pt Arial" font, the same one already
# used for the string/number cells throughout the data rows.
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").Font.ColorIndex = 1
$ws.Range("F1:K1").Font.ColorIndex = -4105

# Move the active selection onto the new first data row, like the edited
# workbook.
$ws.Range("A2:K2").Select()
